$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph (2nd paragraph, right after the
#    title heading) - it is being relocated to the end of the document.
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.Delete()

# 2) Insert a new paragraph just before the final "Prompt: ..." paragraph
#    containing the (bold) title text that used to precede "Meta description".
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
[void]$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 50 Dragons Slot Free - Mesmerizing Design with Lucrative Wins</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newPara.Range.InsertXML($newParaXml)

# 3) Replace the text of the final paragraph's italic run (the old image
#    prompt) with the new meta-description sentence, keeping its formatting.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
[void]$finalRange.MoveEnd(1, -1)
$finalRange.Text = "Experience an immersive adventure - Play 50 Dragons for free and enjoy stunning graphics, bonus features and a potential win of up to €125,000."
